# feat: add 2022-Q3 data
$wb = $excel.ActiveWorkbook

# --- 1. Update the "总计" (summary) sheet: insert a row for 2022-Q3 on top,
#        pushing the existing 2021-Q4 row down to row 3 ---
$wsTotal = $wb.Worksheets.Item("总计")

# Move row 2 (value + style) down to row 3, then overwrite row 2's values.
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))
$wsTotal.Range("A3").Value = 1

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 16
$wsTotal.Range("D2").Value = 1.72

# --- 2. Insert a new worksheet "2022-Q3" right after "总计" (before "2021-Q4") ---
$wsOld = $wb.Worksheets.Item("2021-Q4")
$wsQ3 = $wb.Worksheets.Add($null, $wsTotal)
$wsQ3.Name = "2022-Q3"

# Re-fetch sheet references after the collection changed so the copy below
# targets the freshly materialized sheet correctly.
$wb = $excel.ActiveWorkbook
$wsOld = $wb.Worksheets.Item("2021-Q4")
$wsQ3 = $wb.Worksheets.Item("2022-Q3")

# Copy header row + column-A styling from the existing "2021-Q4" sheet so the
# new sheet matches the workbook's existing look (bold/bordered header, bold
# bordered index column).
$wsOld.Range("B1:H1").Copy($wsQ3.Range("B1:H1"))
$wsOld.Range("A2").Copy($wsQ3.Range("A2:A17"))

$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

# Fund codes / size / position figures are textual in the source data
# (leading zeros, fixed decimal formatting) - force Text format before
# writing so values like "000369" or "91.40" aren't coerced into numbers.
$wsQ3.Range("B2:B17").NumberFormat = "@"
$wsQ3.Range("D2:G17").NumberFormat = "@"

$rows = @(
    @(0, "513500", "博时标普500ETF（QDII）", "71.37", "96.44", "1.51", "1.0777", 8),
    @(1, "000369", "广发全球医疗保健（QDII）人民币A", "2.76", "83.19", "7.16", "0.1976", 1),
    @(2, "000370", "广发全球医疗保健（QDII）美元A", "2.75", "83.19", "7.16", "0.1969", 1),
    @(3, "161125", "易方达标普500指数（QDII-LOF）人民币", "4.74", "90.72", "1.43", "0.0678", 8),
    @(4, "012860", "易方达标普500指数（QDII-LOF）人民币 C", "4.74", "90.72", "1.43", "0.0678", 8),
    @(5, "003718", "易方达标普500指数（QDII-LOF）美元A", "4.66", "90.72", "1.43", "0.0666", 8),
    @(6, "161126", "易方达标普医疗保健指数（QDII-LOF）人民币", "0.57", "93.29", "1.56", "0.0089", 8),
    @(7, "012864", "易方达标普医疗保健指数（QDII-LOF）人民币 C", "0.57", "93.29", "1.56", "0.0089", 8),
    @(8, "159612", "国泰标普500ETF（QDII）", "0.55", "91.40", "1.41", "0.0078", 8),
    @(9, "003719", "易方达标普医疗保健指数（QDII-LOF）美元A", "0.45", "93.29", "1.56", "0.0070", 8),
    @(10, "011706", "长信美国标准普尔100等权重指数增强（QDII）美元", "0.39", "82.64", "0.87", "0.0034", 6),
    @(11, "519981", "长信美国标准普尔100等权重指数增强（QDII）人民币", "0.39", "82.64", "0.87", "0.0034", 6),
    @(12, "012865", "易方达标普医疗保健指数（QDII-LOF）美元 C", "0.12", "93.29", "1.56", "0.0019", 8),
    @(13, "016280", "广发全球医疗保健（QDII）人民币C", "0.02", "83.19", "7.16", "0.0014", 1),
    @(14, "016281", "广发全球医疗保健（QDII）美元C", "0.02", "83.19", "7.16", "0.0014", 1),
    @(15, "012861", "易方达标普500指数（QDII-LOF）美元 C", "0.08", "90.72", "1.43", "0.0011", 8)
)

$r = 2
foreach ($row in $rows) {
    $wsQ3.Range("A$r").Value = $row[0]
    $wsQ3.Range("B$r").Value = $row[1]
    $wsQ3.Range("C$r").Value = $row[2]
    $wsQ3.Range("D$r").Value = $row[3]
    $wsQ3.Range("E$r").Value = $row[4]
    $wsQ3.Range("F$r").Value = $row[5]
    $wsQ3.Range("G$r").Value = $row[6]
    $wsQ3.Range("H$r").Value = $row[7]
    $r++
}
